# Apply updated cryptocurrency price/volume data to Sheet1.
# Mirrors a refreshed scrape: most rows keep their rank, rows 37/38
# (OKB <-> Hedera) and 41/42 (Kaspa <-> dogwifhat) swap rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.643.49"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "3.241.52"
$ws.Range("E3").Value = "  +1.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.15"
$ws.Range("E5").Value = "  +1.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.08"
$ws.Range("E6").Value = "  +2.10%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.243.01"
$ws.Range("E8").Value = "  +1.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  +2.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.65"
$ws.Range("E11").Value = "  -7.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.510"
$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.97"
$ws.Range("E14").Value = "  +0.14%  "

$ws.Range("D15").Value = "3.779.43"
$ws.Range("E15").Value = "  +1.56%  "

$ws.Range("D16").Value = "66.783.65"
$ws.Range("E16").Value = "  +0.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.50"
$ws.Range("E17").Value = "  +0.82%  "

$ws.Range("D18").Value = "3.254.95"
$ws.Range("E18").Value = "  +1.71%  "

$ws.Range("E19").Value = "  +1.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "512.03"
$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.37"
$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.09"
$ws.Range("E23").Value = "  +1.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.91"
$ws.Range("E24").Value = "  -1.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.81"
$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.34"
$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.01"
$ws.Range("E28").Value = "  +0.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.41"
$ws.Range("E29").Value = "  +5.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.04"
$ws.Range("E30").Value = "  +4.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.10"
$ws.Range("E31").Value = "  +1.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.27"
$ws.Range("E32").Value = "  +0.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.18"
$ws.Range("E34").Value = "  -3.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.54"
$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "520.95"
$ws.Range("E36").Value = "  +7.40%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0950"
$ws.Range("E37").Value = "  +5.53%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.16"
$ws.Range("E38").Value = "  +2.48%  "

$ws.Range("D39").Value = "0.0₃0758"
$ws.Range("E39").Value = "  +16.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0420"
$ws.Range("E40").Value = "  +0.55%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.129"
$ws.Range("E41").Value = "  +5.09%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.99"
$ws.Range("E42").Value = "  +2.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.83"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.302"
$ws.Range("E44").Value = "  +1.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  +3.60%  "

$ws.Range("D46").Value = "2.864.37"
$ws.Range("E46").Value = "  -2.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.54"
$ws.Range("E47").Value = "  +0.39%  "

$ws.Range("E48").Value = "  +3.97%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("E51").Value = "  +0.99%  "
